# Edit LOM3046.xlsx worksheet to add new course content rows
# (Objetivos text, two docentes rows split out, Programa resumido/completo
# text, corrected Metodo/Criterio/Norma content, and a new Bibliografia
# paragraph), matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new blank rows at row 13 (pushes old rows 13-23 down to 15-25).
#    This lines up every row height below the insertion point with the target
#    layout, so no further height changes are required for rows 15-25.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Resize(2).Insert()

# The insert operation copies column A's formatting (bold label style) down
# into the two new rows even though they should have no value/style in
# column A. Clear those two cells completely so they disappear from the
# sheet, matching rows 13/14 in the target (B/C only, no A).
$ws.Range("A13:A14").Clear()

# ---------------------------------------------------------------------------
# 2) Give the new B13:C14 cells the same formatting as the other B/C data
#    cells (normal wrap font for B, red wrap font for C) by copying the
#    format from row 15 (an existing fully-formatted data row).
# ---------------------------------------------------------------------------
$ws.Range("B15").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the corrected / new cell values.
# ---------------------------------------------------------------------------
# Objetivos: real objectives paragraph (previously wrongly held a professor's name as a placeholder).
$ws.Range("B10").Value = 'Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural.'
$ws.Range("C10").Value = 'Apresentação introdutória das técnicas de análise microestrutural de materiais. Apresentação das técnicas e equipamentos necessários para a análise microestrutural. Seleção adequada das técnicas experimentais. Apresentação das técnicas adequadas de preparação de amostras. Verificação dos custos envolvidos nas técnicas de caracterização microestrutural.'
# Docentes responsaveis, professor 1 of 2 (now split onto its own row).
$ws.Range("B13").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C13").Value = '6495737 - Durval Rodrigues Junior'
# Docentes responsaveis, professor 2 of 2 (now split onto its own row).
$ws.Range("B14").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("C14").Value = '1643715 - Paulo Atsushi Suzuki'
# Programa resumido: real short-syllabus paragraph (previously held a placeholder date string).
$ws.Range("B15").Value = 'A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais.'
$ws.Range("C15").Value = 'A Microestrutura dos Materiais. Difratometria de raios X. Análise Microestrutural utilizando Luz Síncrotron. Microscopia Óptica. Microscopia Eletrônica. Microscopia de Tunelamento e de Força Atômica. Análise Química de Microrregiões. Análises Térmicas. Fluorescência de raios X. Técnicas Indiretas de Análise de Microestrutura. Seleção de Técnicas Experimentais.'
# Programa: real full syllabus paragraph (previously held a placeholder professor's name).
$ws.Range("B17").Value = '1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais.'
$ws.Range("C17").Value = '1. A Microestrutura dos Materiais. 2. Difratometria de raios X. 3. Análise Microestrutural utilizando Luz Síncrotron. 4. Microscopia Óptica. 5. Microscopia Eletrônica. 6. Microscopia de Tunelamento e de Força Atômica. 7. Análise Química de Microrregiões. 8. Análises Térmicas. 9. Fluorescência de raios X. 10. Técnicas Indiretas de Análise de Microestrutura. 11. Seleção de Técnicas Experimentais.'
# Metodo: real assessment-method paragraph (previously held a placeholder professor's name).
$ws.Range("B20").Value = 'Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre).'
$ws.Range("C20").Value = 'Aplicação de duas avaliações escritas (Aval1 e Aval2) e entrega de relatórios sobre as atividades experimentais. As avaliações e relatórios dividirão o período letivo em dois bimestres. Duas notas (P1 e P2), sendo uma em cada bimestre, serão calculadas como Pn = 0,80 x Avaln + 0,20 x (média aritmética dos relatórios do bimestre).'
# Criterio: real grading-criteria paragraph (previously held the Metodo paragraph).
$ws.Range("B21").Value = 'A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2.'
$ws.Range("C21").Value = 'A Nota Final (NF) do semestre, chamada de primeira avaliação, será a média aritmética das notas P1 e P2.'
# Norma de recuperacao: real make-up-exam paragraph (previously held the Criterio paragraph).
$ws.Range("B22").Value = 'Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação).'
$ws.Range("C22").Value = 'Aplicação de prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final do semestre (primeira avaliação).'
# Bibliografia: real bibliography paragraph (previously held the Norma de recuperacao paragraph).
$ws.Range("B23").Value = '1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. 2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. 3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. 4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. 5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. 7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. 8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. 9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. 10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. 11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. 12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. 13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. 14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. 15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.'
$ws.Range("C23").Value = '1. Van Vlack, L.H. Princípios de Ciência e Tecnologia dos Materiais, 4a.ed., Ed. Campus, Rio de Janeiro, 1984. 2. Shackelford, J.F. Introduction to Materials Science for Engineers. 4th Edition. Prentice Hall Inc., 1996. 3. Padilha, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985. 4. Guy, A.G. Ciência dos Materiais. Livros Técnicos e Científicos Editora, 1982. 5. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 6. Nondestructive Characterization of Materials. Series. Plenum Press, New York. 7. Yacobi, B.G. Holt, D.B. Kazmerski, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994. 8. Lowell, S.; Shields, J. E.; Thomas, M. A.; Thommes, M. Characterization of Porous Solids and Powders: Surface Area, Pore Size and Density, Springer, 2010. 9. Murphy, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001. 10. Wu, Q.; Merchant, F.; Castleman, K. Microscope Image Processing, Academic Press, 2008. 11. Cullity, B. D.; Stock, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001. 12. Goldstein, J.; et al., Scanning Electron Microscopy and X-ray Microanalysis, Springer, 2003. 13. Hatakeyama, T.; Zhenhai, L. Handbook of Thermal Analysis, NY: Wiley, 1999. 14. Haines, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002. 15. Schramm, G. Reologia e Reometria. Editora Artliber, 2006.16. Azevedo, A. D.; Mothe, C. G. Análise Térmica de Materiais. São Paulo: ARTLIBER, 2009.17. Brown, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.18. Muller, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.19. Speyer, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.'
